$d = $word.ActiveDocument

function Insert-XmlReplacingRange($rng, [string]$bodyXml) {
    $xml = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
$bodyXml
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@
    $rng.InsertXML($xml)
}

# ---------------------------------------------------------------------------
# 1) Before the "Design" Heading1 paragraph: add a new "Preventing overflow"
#    Heading2 section (which also takes over the lastRenderedPageBreak that
#    used to sit on the "Design" run) plus its body paragraph. Then strip the
#    lastRenderedPageBreak off the (now second) "Design" run.
# ---------------------------------------------------------------------------

$designPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text.TrimEnd([char]13, [char]7) -eq "Design" -and $cand.Style.NameLocal -eq "Heading 1") {
        $designPara = $cand
        break
    }
}

$designIndex = $designPara.Index
$insertRng = $designPara.Range
$insertRng.Collapse(1)
$insertRng.InsertParagraphBefore()

$placeholder = $d.Paragraphs.Item($designIndex)
$placeholderRng = $placeholder.Range

$overflowBody = @'
<w:p>
<w:pPr>
<w:pStyle w:val="Heading2"/>
<w:rPr>
<w:lang w:val="en-US"/>
</w:rPr>
</w:pPr>
<w:r>
<w:rPr>
<w:lang w:val="en-US"/>
</w:rPr>
<w:lastRenderedPageBreak/>
<w:t>Preventing overflow</w:t>
</w:r>
</w:p>
<w:p>
<w:pPr>
<w:rPr>
<w:lang w:val="en-US"/>
</w:rPr>
</w:pPr>
<w:r>
<w:rPr>
<w:lang w:val="en-US"/>
</w:rPr>
<w:t>The extra pump inside the sump tank is used to manually water plants outside of the system. It will also be used to prevent overflow in the sump tank. The sump tank is connected with an outside system and sends water there when it overflows, but in case the connection is blocked, the extra pump will be used. After the water reaches a level the pump will turn on for some time. It will not turn off when the water drops under the float sensor, because the change will be instantaneous and the pump will turn on/off continuously, risking a damage to the pump.</w:t>
</w:r>
</w:p>
'@

Insert-XmlReplacingRange $placeholderRng $overflowBody

# Remove the lastRenderedPageBreak from the "Design" run (it moved to the
# new "Preventing overflow" paragraph above). Re-locate it by content/style
# rather than by index, since the XML insertion above shifted indices.
$designPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text.TrimEnd([char]13, [char]7) -eq "Design" -and $cand.Style.NameLocal -eq "Heading 1") {
        $designPara = $cand
        break
    }
}
$designRng = $designPara.Range
$designBody = @'
<w:p>
<w:pPr>
<w:pStyle w:val="Heading1"/>
</w:pPr>
<w:r>
<w:t>Design</w:t>
</w:r>
</w:p>
'@
Insert-XmlReplacingRange $designRng $designBody

# ---------------------------------------------------------------------------
# 2) Rework the "Fan" paragraph text (split "FanInterval" -> separate
#    fanIntervalOn / fanIntervalOff run pieces) and append the new P1, P2
#    and NOTE paragraphs, plus a trailing empty paragraph.
# ---------------------------------------------------------------------------

$fanPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -like "Fan:*") {
        $fanPara = $cand
        break
    }
}

$fanRng = $fanPara.Range
$fanBody = @'
<w:p>
<w:r>
<w:rPr>
<w:b/>
</w:rPr>
<w:t>Fan</w:t>
</w:r>
<w:r>
<w:t>: When airTemperature &gt;= airTemperatureThreshold then fan turns on for interval</w:t>
</w:r>
<w:r>
<w:t>s</w:t>
</w:r>
<w:r>
<w:t xml:space="preserve"> (i.e. on for </w:t>
</w:r>
<w:r>
<w:t>f</w:t>
</w:r>
<w:r>
<w:t>anInterval</w:t>
</w:r>
<w:r>
<w:t>On</w:t>
</w:r>
<w:r>
<w:t xml:space="preserve"> seconds and off for </w:t>
</w:r>
<w:r>
<w:t>f</w:t>
</w:r>
<w:r>
<w:t>anInterval</w:t>
</w:r>
<w:r>
<w:t>Off</w:t>
</w:r>
<w:r>
<w:t xml:space="preserve"> seconds), until airTemperature &lt; airTemperatureThreshold.</w:t>
</w:r>
</w:p>
<w:p>
<w:r>
<w:rPr>
<w:b/>
</w:rPr>
<w:t>P1</w:t>
</w:r>
<w:r>
<w:rPr>
<w:b/>
</w:rPr>
<w:t>:</w:t>
</w:r>
<w:r>
<w:t xml:space="preserve"> </w:t>
</w:r>
<w:r>
<w:t>The pump will turn off regardless of temperature when FloatLow = off (</w:t>
</w:r>
<w:r>
<w:rPr>
<w:b/>
<w:i/>
</w:rPr>
<w:t>highest priority</w:t>
</w:r>
<w:r>
<w:t xml:space="preserve">). </w:t>
</w:r>
<w:r>
<w:t xml:space="preserve">When </w:t>
</w:r>
<w:r>
<w:t>waterTemperature &lt;= waterTemperatureThreshold then turn pump on for intervals (i.e. on for p1IntervalOn seconds and off for p1IntervalOff seconds) until waterTemperature &gt; waterTemperatureThreshold.</w:t>
</w:r>
</w:p>
<w:p>
<w:r>
<w:rPr>
<w:b/>
</w:rPr>
<w:t>P2</w:t>
</w:r>
<w:r>
<w:t>: When FloatHigh = on turn on for p2IntervalOn and then turn off.</w:t>
</w:r>
</w:p>
<w:p>
<w:pPr>
<w:rPr>
<w:b/>
</w:rPr>
</w:pPr>
<w:r>
<w:rPr>
<w:b/>
</w:rPr>
<w:t>NOTE: Pumps and fan will have override switches to be manually turned on.</w:t>
</w:r>
</w:p>
<w:p/>
'@

Insert-XmlReplacingRange $fanRng $fanBody
